# Gitsunum.pptx - slide 2 ("Git") body placeholder:
# Fix/expand the sentence about the "checkpoint" analogy:
#   " sistemini andıran Git temel olarak üzerinde çalıştığımız projeyi
#     adım adım izlememize yardım olur."
# becomes
#   " sistemini andıran Git, temel olarak üzerinde çalıştığımız projeyi
#     adım adım izlememize yardımcı olur."
# (adds a comma after "Git" and turns "yardım" into "yardımcı").

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(2)
$shp = $s.Shapes.Item(2)
$tr  = $shp.TextFrame.TextRange

$oldSentence = " sistemini andıran Git temel olarak üzerinde çalıştığımız projeyi adım adım izlememize yardım olur."
$newSentence = " sistemini andıran Git, temel olarak üzerinde çalıştığımız projeyi adım adım izlememize yardımcı olur."

$fullText = $tr.Text
$startIdx = $fullText.IndexOf($oldSentence)

if ($startIdx -ge 0) {
    # Characters() is 1-based; target only the run that holds this sentence
    # so the rest of the paragraph/shape is left completely untouched.
    $run = $tr.Characters($startIdx + 1, $oldSentence.Length)
    $run.Text = $newSentence
}
